# The workbook tracks daily Tuna prices for "Terminal La Palmera de La Serena".
# Two new daily records need to be inserted before the current row 36, which
# pushes the existing rows 36-44 down to 38-46 (dimension grows from T44 to T46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above the current row 36 (shifts 36..44 -> 38..46).
$ws.Rows.Item(36).EntireRow.Insert()
$ws.Rows.Item(36).EntireRow.Insert()

# New row 36: Especial quality record dated 45005.
$ws.Range("A36").Value = 8
$ws.Range("B36").Value = "Terminal La Palmera de La Serena"
$ws.Range("C36").Value = "Coquimbo"
$ws.Range("D36").Value = 45005
$ws.Range("E36").Value = 4
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100107
$ws.Range("H36").Value = "Otros"
$ws.Range("I36").Value = 100107011
$ws.Range("J36").Value = "Tuna"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Especial"
$ws.Range("M36").Value = 300
$ws.Range("N36").Value = 13000
$ws.Range("O36").Value = 14000
$ws.Range("P36").Value = 13500
$ws.Range("Q36").Value = "$/caja 18 kilos"
$ws.Range("R36").Value = "Provincia de Limarí"
$ws.Range("S36").Value = 750
$ws.Range("T36").Value = 18

# New row 37: Primera quality record dated 45005.
$ws.Range("A37").Value = 8
$ws.Range("B37").Value = "Terminal La Palmera de La Serena"
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = 45005
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100107
$ws.Range("H37").Value = "Otros"
$ws.Range("I37").Value = 100107011
$ws.Range("J37").Value = "Tuna"
$ws.Range("K37").Value = "Sin especificar"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 200
$ws.Range("N37").Value = 11000
$ws.Range("O37").Value = 12000
$ws.Range("P37").Value = 11500
$ws.Range("Q37").Value = "$/caja 18 kilos"
$ws.Range("R37").Value = "Provincia de Limarí"
$ws.Range("S37").Value = 639
$ws.Range("T37").Value = 18
